$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6534678339958191
$ws.Range("B1").Value = 0.9482837915420532
$ws.Range("C1").Value = 4.09104061126709
$ws.Range("D1").Value = 2.124695301055908
$ws.Range("E1").Value = 1.654451727867126
